$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") entirely; remaining rows shift up by one.
$ws.Rows("2").Delete()
